$d = $word.ActiveDocument

$replacements = @(
    @("25×54=", "33×20="),
    @("32×64=", "24×44="),
    @("18×34=", "65×62="),
    @("99×61=", "40×53="),
    @("93×51=", "92×98="),
    @("23×69=", "53×64="),
    @("49×40=", "38×93="),
    @("98×23=", "84×79="),
    @("56×77=", "99×39="),
    @("85×99=", "75×64="),
    @("12×15=", "75×96="),
    @("60×28=", "59×74="),
    @("25×60=", "54×99="),
    @("20×18=", "26×35="),
    @("75×26=", "54×41="),
    @("60×94=", "84×24="),
    @("20×99=", "30×73="),
    @("86×48=", "49×93="),
    @("63×37=", "62×54="),
    @("47×30=", "40×87="),
    @("19×39=", "36×65="),
    @("58×32=", "43×57="),
    @("52×15=", "18×55="),
    @("46×71=", "18×93="),
    @("25×21=", "28×20=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
